$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLUESKY-COMMAND-TABLE")

# ---------------------------------------------------------------------------
# Synonyms table (rows 88-115): add a "Description" column (C) that shows the
# description of the equivalent command already referenced in column B, plus
# an (empty) column D so the table is a full 4-column block like the rest of
# the sheet. Rows 97/98 ("END"/"EXIT") pointed at the now-removed "STOP"
# synonym row and are repointed at "QUIT" (their real equivalent command).
# ---------------------------------------------------------------------------

# Fix the two rows whose equivalent command changes from STOP to QUIT.
$ws.Range("B97").Value = "QUIT"
$ws.Range("B98").Value = "QUIT"

# --- Header row (row 88): "Description" label in C88, blank styled D88 -----
$ws.Range("C88").Value = "Description"

$ws.Range("A88").Copy()
$ws.Range("C88").PasteSpecial(-4122)
$ws.Range("C88").Value = "Description"
$ws.Range("C88").Font.Bold = $false
$ws.Range("C88").HorizontalAlignment = -4131
$ws.Range("C88").Borders.LineStyle = -4142

$ws.Range("C88").Copy()
$ws.Range("D88").PasteSpecial(-4122)
$ws.Range("D88").Value = ""

# --- Data rows (89-115): Description of the equivalent command in col B ---
$descriptions = @{
    89  = "Show help in a command or write list of commands to file"
    90  = "Quit program/Stop simulation"
    91  = "Start/Run simulation or continue after pause"
    92  = "Create an aircraft"
    93  = "Delete command (aircraft, wind, area)"
    94  = "Go direct to specified waypoint in route (FMS)"
    95  = "Go direct to specified waypoint in route (FMS)"
    96  = "Switch on/off elements and background of map/radar view"
    97  = "Quit program/Stop simulation"
    98  = "Quit program/Stop simulation"
    99  = "Fast forward the simulation"
    100 = "Set resolution method to be used horizontally"
    101 = "Set resolution method to be used horizontally"
    102 = "Set resolution method to be used horizontally"
    103 = "Initial condition: (re)start simulation and open scenario file"
    104 = "Initial condition: (re)start simulation and open scenario file"
    105 = "Pause(hold) simulation"
    106 = "Quit program/Stop simulation"
    107 = "Set resolution factor horizontal (to add a margin)"
    108 = "Set resolution factor vertical (to add a margin)"
    109 = "Start/Run simulation or continue after pause"
    110 = "Start/Run simulation or continue after pause"
    111 = "Quit program/Stop simulation"
    112 = "Heading command (autopilot)"
    113 = "Set resolution method to be used vertically"
    114 = "Set resolution method to be used vertically"
    115 = "Set resolution method to be used vertically"
}

foreach ($row in 89..115) {
    $aCell = "A$row"
    $cCell = "C$row"
    $dCell = "D$row"

    $ws.Range($aCell).Copy()
    $ws.Range($cCell).PasteSpecial(-4122)
    $ws.Range($cCell).Value = $descriptions[$row]
    $ws.Range($cCell).Font.Bold = $false
    $ws.Range($cCell).WrapText = $false

    $ws.Range("B$row").Copy()
    $ws.Range($dCell).PasteSpecial(-4122)
    $ws.Range($dCell).Value = ""
}

# --- Match the saved sheet view (scroll position / active selection) ------
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("F92").Select()
